$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$count = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $oldVal) {
        $cell.Value2 = $newVal
        $count++
    }
}

Write-Host "Replaced $count cells"
